$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reverse the order of the "Periodo Mora" period labels (E16:E22) so they
# run from the most recent period (2311) down to the oldest (2305).
$ws.Range("E16").Value = "2311"
$ws.Range("E17").Value = "2310"
$ws.Range("E18").Value = "2309"
$ws.Range("E19").Value = "2308"
$ws.Range("E20").Value = "2307"
$ws.Range("E21").Value = "2306"
$ws.Range("E22").Value = "2305"

# Swap the corresponding dates in column F to match the reordered periods.
$ws.Range("F16").Value = 43307
$ws.Range("F22").Value = 40214
